$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Q0) ---
$ws.Range("B2").Value = 0.2048957379308303
$ws.Range("C2").Value = 0.3997631070090498
$ws.Range("D2").Value = 0.3041926803226083
$ws.Range("E2").Value = 0.5515366536528722
$ws.Range("F2").Value = 0.5313948145867153

# --- Row 3 (Q1) ---
$ws.Range("B3").Value = 0.2376769736950485
$ws.Range("C3").Value = 0.3882757705769553
$ws.Range("D3").Value = 0.3054694716933793
$ws.Range("E3").Value = 0.552692927124438
$ws.Range("F3").Value = 0.5193528555079431

# --- Row 4 (Q2) ---
$ws.Range("B4").Value = 0.2742536187534859
$ws.Range("C4").Value = 0.4315080895772995
$ws.Range("D4").Value = 0.2558436790103293
$ws.Range("E4").Value = 0.5058099237958161
$ws.Range("F4").Value = 0.4439024851280417

# --- Row 5 (Q3) ---
$ws.Range("B5").Value = 0.3352241573957673
$ws.Range("C5").Value = 0.3891781958695149
$ws.Range("D5").Value = 0.2638597261285222
$ws.Range("E5").Value = 0.5136727811832376
$ws.Range("F5").Value = 0.4082069811621331

# --- Row 6 (Q4) ---
$ws.Range("B6").Value = 0.3275244721553913
$ws.Range("C6").Value = 0.4092522750510487
$ws.Range("D6").Value = 0.3133210581871232
$ws.Range("E6").Value = 0.5597508894027086
$ws.Range("F6").Value = 0.4784799755782838

# --- Row 7 (Q5) ---
$ws.Range("B7").Value = 0.3138946273557288
$ws.Range("C7").Value = 0.3664898380679329
$ws.Range("D7").Value = 0.2144244288482471
$ws.Range("E7").Value = 0.4630598544985812
$ws.Range("F7").Value = 0.3610836686090043
$ws.Range("G7").Value = 9

# --- Row 8 (Q6) ---
$ws.Range("B8").Value = 0.3484515249888549
$ws.Range("C8").Value = 0.349600236492926
$ws.Range("D8").Value = 0.3199050252412782
$ws.Range("E8").Value = 0.5656014720996386
$ws.Range("F8").Value = 0.4880408507175025
$ws.Range("G8").Value = 6

# --- Row 9 (Q7) ---
$ws.Range("B9").Value = 0.6949929595738982
$ws.Range("C9").Value = 0.6949929595738982
$ws.Range("D9").Value = 0.5900249945076564
$ws.Range("E9").Value = 0.7681308446532117
$ws.Range("F9").Value = 0.4006428221939781
$ws.Range("G9").Value = 3

# --- Row 10 (new, Q8) ---
# Copy formatting from the row above (A9, which carries the bold/border/centered style)
# so that the new label cell A10 matches the style used by the other row labels.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = 0.1275310031338272
$ws.Range("C10").Value = 0.1275310031338272
$ws.Range("D10").Value = 0.01626415676032024
$ws.Range("E10").Value = 0.1275310031338272
$ws.Range("G10").Value = 1
